$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.384.08"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.012.01"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.78"
$ws.Range("E5").Value = "  +4.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.16"
$ws.Range("E8").Value = "  -6.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0769"
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "2.308.76"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  -5.17%  "
$ws.Range("E14").Value = "  -4.61%  "
$ws.Range("E15").Value = "  -7.74%  "
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "2.016.75"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "37.236.25"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "0.0₃0838"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.09"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.64"
$ws.Range("E23").Value = "  +6.77%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.36"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.67"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -9.00%  "
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.82"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0931"
$ws.Range("E42").Value = "  -5.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "1.390.21"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.98"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.67"
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  -4.54%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "2.201.57"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("E51").Value = "  -4.03%  "
